# Apply the "harvard case classification" edit:
# - Insert a new "average_doctor_old" column at BP (shifting the former
#   "average_doctor" data into BQ, now labeled "average_doctor"),
# - Recompute/update numeric stats across the data rows (4-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the BP/BQ column labels ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Data rows: updated statistics ---

# Row 4
$ws.Range("E4").Value = 0.445
$ws.Range("F4").Value = 0.067
$ws.Range("G4").Value = 0.259
$ws.Range("N4").Value = 0.438
$ws.Range("O4").Value = 0.065
$ws.Range("P4").Value = 0.256
$ws.Range("Q4").Value = 0.021
$ws.Range("R4").Value = 0.015
$ws.Range("S4").Value = 0.124
$ws.Range("W4").Value = 0.305
$ws.Range("X4").Value = 0.113
$ws.Range("Y4").Value = 0.336
$ws.Range("AI4").Value = 0.36
$ws.Range("AJ4").Value = 0.08599999999999999
$ws.Range("AK4").Value = 0.293
$ws.Range("AU4").Value = 0.207
$ws.Range("AV4").Value = 0.031
$ws.Range("AW4").Value = 0.175
$ws.Range("BA4").Value = 1.963
$ws.Range("BB4").Value = 0.154
$ws.Range("BC4").Value = 0.392
$ws.Range("BG4").Value = 0.717
$ws.Range("BH4").Value = 0.135
$ws.Range("BI4").Value = 0.368
$ws.Range("BM4").Value = 0.709
$ws.Range("BN4").Value = 0.075
$ws.Range("BO4").Value = 0.274
$ws.Range("BP4").Value = 0.654
$ws.Range("BQ4").Value = 0.6879999999999999

# Row 5
$ws.Range("E5").Value = 0.5620000000000001
$ws.Range("F5").Value = 0.07099999999999999
$ws.Range("G5").Value = 0.267
$ws.Range("N5").Value = 0.732
$ws.Range("O5").Value = 0.08
$ws.Range("P5").Value = 0.283
$ws.Range("Q5").Value = 0.011
$ws.Range("R5").Value = 0.003
$ws.Range("S5").Value = 0.053
$ws.Range("W5").Value = 0.287
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.33
$ws.Range("AI5").Value = 0.384
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.309
$ws.Range("AU5").Value = 0.394
$ws.Range("AV5").Value = 0.1
$ws.Range("AW5").Value = 0.316
$ws.Range("BA5").Value = 1.302
$ws.Range("BB5").Value = 0.079
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.383
$ws.Range("BH5").Value = 0.043
$ws.Range("BI5").Value = 0.207
$ws.Range("BM5").Value = 0.541
$ws.Range("BN5").Value = 0.061
$ws.Range("BO5").Value = 0.246
$ws.Range("BP5").Value = 0.434
$ws.Range("BQ5").Value = 0.444

# Row 6
$ws.Range("E6").Value = 0.497
$ws.Range("N6").Value = 0.548
$ws.Range("Q6").Value = 0.014
$ws.Range("W6").Value = 0.296
$ws.Range("AI6").Value = 0.372
$ws.Range("AU6").Value = 0.271
$ws.Range("BA6").Value = 1.557
$ws.Range("BG6").Value = 0.499
$ws.Range("BM6").Value = 0.614
$ws.Range("BP6").Value = 0.519
$ws.Range("BQ6").Value = 0.537

# Row 7
$ws.Range("E7").Value = 0.534
$ws.Range("N7").Value = 0.645
$ws.Range("Q7").Value = 0.012
$ws.Range("W7").Value = 0.29
$ws.Range("AI7").Value = 0.379
$ws.Range("AU7").Value = 0.334
$ws.Range("BA7").Value = 1.392
$ws.Range("BG7").Value = 0.422
$ws.Range("BM7").Value = 0.5679999999999999
$ws.Range("BP7").Value = 0.464
$ws.Range("BQ7").Value = 0.476

# Row 8
$ws.Range("E8").Value = 0.624
$ws.Range("F8").Value = 0.094
$ws.Range("G8").Value = 0.307
$ws.Range("N8").Value = 0.776
$ws.Range("O8").Value = 0.064
$ws.Range("P8").Value = 0.254
$ws.Range("Q8").Value = 0.011
$ws.Range("S8").Value = 0.08
$ws.Range("W8").Value = 0.332
$ws.Range("X8").Value = 0.125
$ws.Range("Y8").Value = 0.354
$ws.Range("AI8").Value = 0.403
$ws.Range("AJ8").Value = 0.125
$ws.Range("AK8").Value = 0.354
$ws.Range("AU8").Value = 0.337
$ws.Range("AW8").Value = 0.296
$ws.Range("BA8").Value = 1.679
$ws.Range("BB8").Value = 0.123
$ws.Range("BC8").Value = 0.35
$ws.Range("BG8").Value = 0.537
$ws.Range("BH8").Value = 0.1
$ws.Range("BI8").Value = 0.317
$ws.Range("BM8").Value = 0.68
$ws.Range("BN8").Value = 0.064
$ws.Range("BO8").Value = 0.252
$ws.Range("BP8").Value = 0.5600000000000001
$ws.Range("BQ8").Value = 0.583

# Row 9
$ws.Range("E9").Value = 0.548
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.671
$ws.Range("O9").Value = 0.221
$ws.Range("P9").Value = 0.47
$ws.Range("W9").Value = 0.233
$ws.Range("X9").Value = 0.179
$ws.Range("Y9").Value = 0.423
$ws.Range("AI9").Value = 0.301
$ws.Range("AJ9").Value = 0.211
$ws.Range("AK9").Value = 0.459
$ws.Range("BA9").Value = 1.575
$ws.Range("BB9").Value = 0.242
$ws.Range("BC9").Value = 0.492
$ws.Range("BG9").Value = 0.548
$ws.Range("BH9").Value = 0.248
$ws.Range("BI9").Value = 0.498
$ws.Range("BM9").Value = 0.616
$ws.Range("BN9").Value = 0.236
$ws.Range("BO9").Value = 0.486
$ws.Range("BP9").Value = 0.525
$ws.Range("BQ9").Value = 0.541

# Row 10
$ws.Range("E10").Value = 0.699
$ws.Range("F10").Value = 0.211
$ws.Range("G10").Value = 0.459
$ws.Range("N10").Value = 0.89
$ws.Range("O10").Value = 0.098
$ws.Range("P10").Value = 0.312
$ws.Range("W10").Value = 0.411
$ws.Range("X10").Value = 0.242
$ws.Range("Y10").Value = 0.492
$ws.Range("AI10").Value = 0.438
$ws.Range("AJ10").Value = 0.246
$ws.Range("AK10").Value = 0.496
$ws.Range("AU10").Value = 0.329
$ws.Range("AV10").Value = 0.221
$ws.Range("AW10").Value = 0.47
$ws.Range("BA10").Value = 2.013
$ws.Range("BB10").Value = 0.248
$ws.Range("BC10").Value = 0.498
$ws.Range("BG10").Value = 0.616
$ws.Range("BH10").Value = 0.236
$ws.Range("BI10").Value = 0.486
$ws.Range("BM10").Value = 0.849
$ws.Range("BN10").Value = 0.128
$ws.Range("BO10").Value = 0.358
$ws.Range("BP10").Value = 0.671
$ws.Range("BQ10").Value = 0.7

# Row 11
$ws.Range("E11").Value = 0.74
$ws.Range("F11").Value = 0.193
$ws.Range("G11").Value = 0.439
$ws.Range("N11").Value = 0.904
$ws.Range("O11").Value = 0.08699999999999999
$ws.Range("P11").Value = 0.294
$ws.Range("W11").Value = 0.411
$ws.Range("X11").Value = 0.242
$ws.Range("Y11").Value = 0.492
$ws.Range("AI11").Value = 0.479
$ws.Range("AJ11").Value = 0.25
$ws.Range("AK11").Value = 0.5
$ws.Range("AU11").Value = 0.452
$ws.Range("AV11").Value = 0.248
$ws.Range("AW11").Value = 0.498
$ws.Range("BA11").Value = 2.013
$ws.Range("BB11").Value = 0.248
$ws.Range("BC11").Value = 0.498
$ws.Range("BG11").Value = 0.616
$ws.Range("BH11").Value = 0.236
$ws.Range("BI11").Value = 0.486
$ws.Range("BM11").Value = 0.849
$ws.Range("BN11").Value = 0.128
$ws.Range("BO11").Value = 0.358
$ws.Range("BP11").Value = 0.671
$ws.Range("BQ11").Value = 0.704

# Row 12
$ws.Range("E12").Value = 1.481
$ws.Range("F12").Value = 0.879
$ws.Range("G12").Value = 0.9379999999999999
$ws.Range("N12").Value = 1.433
$ws.Range("O12").Value = 0.783
$ws.Range("P12").Value = 0.885
$ws.Range("W12").Value = 1.6
$ws.Range("X12").Value = 0.573
$ws.Range("Y12").Value = 0.757
$ws.Range("AI12").Value = 1.743
$ws.Range("AJ12").Value = 1.334
$ws.Range("AK12").Value = 1.155
$ws.Range("AU12").Value = 2.714
$ws.Range("AV12").Value = 2.947
$ws.Range("AW12").Value = 1.717
$ws.Range("BA12").Value = 3.838
$ws.Range("BB12").Value = 0.428
$ws.Range("BC12").Value = 0.654
$ws.Range("BG12").Value = 1.133
$ws.Range("BH12").Value = 0.16
$ws.Range("BI12").Value = 0.4
$ws.Range("BM12").Value = 1.355
$ws.Range("BN12").Value = 0.39
$ws.Range("BO12").Value = 0.625
$ws.Range("BP12").Value = 1.279
$ws.Range("BQ12").Value = 1.303

# Row 13
$ws.Range("E13").Value = 1.547
$ws.Range("F13").Value = 0.635
$ws.Range("G13").Value = 0.797
$ws.Range("N13").Value = 2.001
$ws.Range("O13").Value = 0.788
$ws.Range("P13").Value = 0.887
$ws.Range("W13").Value = 1.026
$ws.Range("X13").Value = 0.194
$ws.Range("Y13").Value = 0.441
$ws.Range("AI13").Value = 1.265
$ws.Range("AJ13").Value = 0.385
$ws.Range("AK13").Value = 0.62
$ws.Range("AU13").Value = 2.172
$ws.Range("AV13").Value = 0.619
$ws.Range("AW13").Value = 0.787
$ws.Range("BA13").Value = 2.326
$ws.Range("BB13").Value = 0.289
$ws.Range("BC13").Value = 0.537
$ws.Range("BG13").Value = 0.584
$ws.Range("BH13").Value = 0.077
$ws.Range("BI13").Value = 0.277
$ws.Range("BM13").Value = 0.887
$ws.Range("BN13").Value = 0.285
$ws.Range("BO13").Value = 0.534
$ws.Range("BP13").Value = 0.775
$ws.Range("BQ13").Value = 0.718
